$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5,D6,D8,D9,D10,D11,D12,D15,D17,D19,D21,D22,D23,D24,D25,D26,D27,D28,D30,D31,D32,D33,D34,D35,D36,D38,D39,D41,D44,D47,D48,D49,D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.536.63"
$ws.Range("E2").Value = "  +5.44%  "

$ws.Range("D3").Value = "1.725.75"
$ws.Range("E3").Value = "  +4.45%  "

$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "225.56"
$ws.Range("E5").Value = "  +3.21%  "

$ws.Range("D6").Value = "0.5369"
$ws.Range("E6").Value = "  +2.94%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "0.2672"
$ws.Range("E8").Value = "  +0.97%  "

$ws.Range("D9").Value = "0.06606"
$ws.Range("E9").Value = "  +4.12%  "

$ws.Range("D10").Value = "21.78"
$ws.Range("E10").Value = "  +6.58%  "

$ws.Range("D11").Value = "0.07720"
$ws.Range("E11").Value = "  +0.32%  "

$ws.Range("D12").Value = "4.610"
$ws.Range("E12").Value = "  -0.39%  "

$ws.Range("D13").Value = "1.729.45"
$ws.Range("E13").Value = "  +3.09%  "

$ws.Range("D14").Value = "1.964.74"
$ws.Range("E14").Value = "  +4.47%  "

$ws.Range("D15").Value = "0.5847"
$ws.Range("E15").Value = "  +4.47%  "

$ws.Range("D16").Value = "0.0₅8304"
$ws.Range("E16").Value = "  +1.65%  "

$ws.Range("D17").Value = "68.02"
$ws.Range("E17").Value = "  +3.93%  "

$ws.Range("D18").Value = "27.561.69"
$ws.Range("E18").Value = "  +5.54%  "

$ws.Range("D19").Value = "221.66"
$ws.Range("E19").Value = "  +15.72%  "

$ws.Range("E20").Value = "  +0.04%  "

$ws.Range("D21").Value = "4.733"
$ws.Range("E21").Value = "  +2.06%  "

$ws.Range("D22").Value = "10.65"
$ws.Range("E22").Value = "  +1.57%  "

$ws.Range("D23").Value = "6.095"
$ws.Range("E23").Value = "  +2.63%  "

$ws.Range("D24").Value = "1.005"
$ws.Range("E24").Value = "  +0.04%  "

$ws.Range("D25").Value = "148.41"
$ws.Range("E25").Value = "  +1.97%  "

$ws.Range("D26").Value = "1.714"
$ws.Range("E26").Value = "  +13.40%  "

$ws.Range("D27").Value = "0.1233"
$ws.Range("E27").Value = "  +3.27%  "

$ws.Range("D28").Value = "7.407"
$ws.Range("E28").Value = "  +2.39%  "

$ws.Range("E29").Value = "  +4.84%  "

$ws.Range("D30").Value = "0.05576"
$ws.Range("E30").Value = "  +1.77%  "

$ws.Range("D31").Value = "1.301"
$ws.Range("E31").Value = "  +2.41%  "

$ws.Range("D32").Value = "3.551"
$ws.Range("E32").Value = "  +2.82%  "

$ws.Range("D33").Value = "3.456"
$ws.Range("E33").Value = "  +2.52%  "

$ws.Range("D34").Value = "1.660"
$ws.Range("E34").Value = "  +6.35%  "

$ws.Range("D35").Value = "0.9637"
$ws.Range("E35").Value = "  +1.17%  "

$ws.Range("D36").Value = "2.828"
$ws.Range("E36").Value = "  +1.48%  "

$ws.Range("E37").Value = "  +1.95%  "

$ws.Range("D38").Value = "0.5949"
$ws.Range("E38").Value = "  +5.50%  "

$ws.Range("D39").Value = "0.01648"
$ws.Range("E39").Value = "  +4.41%  "

$ws.Range("E40").Value = "  +1.07%  "

$ws.Range("D41").Value = "0.8588"
$ws.Range("E41").Value = "  +2.97%  "

$ws.Range("D42").Value = "1.055.02"
$ws.Range("E42").Value = "  +2.45%  "

$ws.Range("E43").Value = "  +0.10%  "

$ws.Range("D44").Value = "101.47"
$ws.Range("E44").Value = "  +0.31%  "

$ws.Range("D45").Value = "1.869.61"
$ws.Range("E45").Value = "  +4.31%  "

$ws.Range("D46").Value = "0.0₈115"
$ws.Range("E46").Value = "  +6.70%  "

$ws.Range("D47").Value = "59.10"
$ws.Range("E47").Value = "  +2.39%  "

$ws.Range("D48").Value = "8.233"
$ws.Range("E48").Value = "  +3.20%  "

$ws.Range("D49").Value = "0.4437"
$ws.Range("E49").Value = "  +2.23%  "

$ws.Range("E50").Value = "  +0.10%  "

$ws.Range("D51").Value = "0.05264"
$ws.Range("E51").Value = "  +1.47%  "

